$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name / Link columns (plain text, safe for direct assignment)
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'

# Update Price / Volume columns via text-literal formulas to avoid numeric auto-conversion
$ws.Range('D2').Formula = '="30.762.68"'
$ws.Range('E2').Formula = '="  +1.58%  "'
$ws.Range('D3').Formula = '="2.107.37"'
$ws.Range('E3').Formula = '="  +5.11%  "'
$ws.Range('D4').Formula = '="1.005"'
$ws.Range('E4').Formula = '="  +0.50%  "'
$ws.Range('D5').Formula = '="333.66"'
$ws.Range('E5').Formula = '="  +2.76%  "'
$ws.Range('D6').Formula = '="1.004"'
$ws.Range('E6').Formula = '="  +0.46%  "'
$ws.Range('D7').Formula = '="0.5292"'
$ws.Range('E7').Formula = '="  +3.73%  "'
$ws.Range('D8').Formula = '="0.4360"'
$ws.Range('E8').Formula = '="  +4.99%  "'
$ws.Range('D9').Formula = '="0.08913"'
$ws.Range('E9').Formula = '="  +2.33%  "'
$ws.Range('D10').Formula = '="47.01"'
$ws.Range('E10').Formula = '="  +9.77%  "'
$ws.Range('D11').Formula = '="1.165"'
$ws.Range('E11').Formula = '="  +2.55%  "'
$ws.Range('D12').Formula = '="24.75"'
$ws.Range('E12').Formula = '="  -1.23%  "'
$ws.Range('D13').Formula = '="2.117.40"'
$ws.Range('E13').Formula = '="  +5.98%  "'
$ws.Range('D14').Formula = '="6.709"'
$ws.Range('E14').Formula = '="  +2.72%  "'
$ws.Range('D15').Formula = '="7.747"'
$ws.Range('E15').Formula = '="  +4.15%  "'
$ws.Range('D16').Formula = '="96.83"'
$ws.Range('E16').Formula = '="  +2.83%  "'
$ws.Range('D17').Formula = '="1.005"'
$ws.Range('E17').Formula = '="  +0.61%  "'
$ws.Range('D18').Formula = '="0.00001130"'
$ws.Range('E18').Formula = '="  +1.21%  "'
$ws.Range('D19').Formula = '="0.06680"'
$ws.Range('E19').Formula = '="  +2.10%  "'
$ws.Range('D20').Formula = '="18.98"'
$ws.Range('E20').Formula = '="  +0.43%  "'
$ws.Range('D21').Formula = '="1.002"'
$ws.Range('E21').Formula = '="  +0.25%  "'
$ws.Range('D22').Formula = '="6.298"'
$ws.Range('E22').Formula = '="  +2.41%  "'
$ws.Range('D23').Formula = '="30.828.95"'
$ws.Range('E23').Formula = '="  +1.59%  "'
$ws.Range('D24').Formula = '="12.18"'
$ws.Range('E24').Formula = '="  +4.01%  "'
$ws.Range('D25').Formula = '="2.308"'
$ws.Range('E25').Formula = '="  +4.31%  "'
$ws.Range('D26').Formula = '="2.365.85"'
$ws.Range('E26').Formula = '="  +6.00%  "'
$ws.Range('D27').Formula = '="22.63"'
$ws.Range('E27').Formula = '="  -0.06%  "'
$ws.Range('D28').Formula = '="2.576"'
$ws.Range('E28').Formula = '="  +7.30%  "'
$ws.Range('D29').Formula = '="162.96"'
$ws.Range('E29').Formula = '="  -0.24%  "'
$ws.Range('D30').Formula = '="132.98"'
$ws.Range('E30').Formula = '="  +1.15%  "'
$ws.Range('D31').Formula = '="1.180"'
$ws.Range('E31').Formula = '="  +3.53%  "'
$ws.Range('D32').Formula = '="0.1078"'
$ws.Range('E32').Formula = '="  +2.32%  "'
$ws.Range('D33').Formula = '="4.055"'
$ws.Range('E33').Formula = '="  +6.35%  "'
$ws.Range('D34').Formula = '="6.174"'
$ws.Range('E34').Formula = '="  +1.67%  "'
$ws.Range('D35').Formula = '="1.535"'
$ws.Range('E35').Formula = '="  +13.87%  "'
$ws.Range('D36').Formula = '="0.02588"'
$ws.Range('E36').Formula = '="  +3.45%  "'
$ws.Range('D37').Formula = '="0.06741"'
$ws.Range('E37').Formula = '="  +2.82%  "'
$ws.Range('D38').Formula = '="9.528"'
$ws.Range('E38').Formula = '="  +6.86%  "'
$ws.Range('D39').Formula = '="5.499"'
$ws.Range('E39').Formula = '="  +1.76%  "'
$ws.Range('D40').Formula = '="0.2274"'
$ws.Range('E40').Formula = '="  +3.52%  "'
$ws.Range('D41').Formula = '="12.60"'
$ws.Range('E41').Formula = '="  +6.27%  "'
$ws.Range('D42').Formula = '="0.6801"'
$ws.Range('E42').Formula = '="  +2.97%  "'
$ws.Range('D43').Formula = '="1.244"'
$ws.Range('E43').Formula = '="  +1.46%  "'
$ws.Range('D44').Formula = '="1.003"'
$ws.Range('E44').Formula = '="  +0.33%  "'
$ws.Range('D45').Formula = '="14.09"'
$ws.Range('E45').Formula = '="  +3.19%  "'
$ws.Range('D46').Formula = '="0.6398"'
$ws.Range('E46').Formula = '="  +4.24%  "'
$ws.Range('D47').Formula = '="2.219"'
$ws.Range('E47').Formula = '="  +0.35%  "'
$ws.Range('D48').Formula = '="3.660"'
$ws.Range('E48').Formula = '="  -0.17%  "'
$ws.Range('D49').Formula = '="1.259"'
$ws.Range('E49').Formula = '="  -0.40%  "'
$ws.Range('D50').Formula = '="1.199"'
$ws.Range('E50').Formula = '="  +9.93%  "'
$ws.Range('D51').Formula = '="82.56"'
$ws.Range('E51').Formula = '="  +3.18%  "'

# Convert the formulas to static values (preserve text type, avoid style churn)
$rng = $ws.Range("D2:E51")
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
